$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.530133605003357
$ws.Range("B1").Value = 1.620538830757141
$ws.Range("C1").Value = 2.183976173400879
$ws.Range("D1").Value = 4.027495861053467
$ws.Range("E1").Value = 2.63889741897583
